# Adds a new "2022-Q3" quarter sheet (placed right after "总计" and before
# "2022-Q2"), populated with the new quarter's per-fund holding data, and
# records the new quarter in the "总计" summary sheet as its first data row
# (shifting the existing summary rows down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for "2022-Q3".
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Copy the formatting of the (now shifted) row 3 onto the newly inserted
# row 2 so the new row matches the existing look (bold/bordered A column,
# plain B/C/D columns) instead of whatever Excel guessed on insert.
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.05

# Re-sequence the leading index column (0-based row counter) for every
# data row now that there are 8 of them (rows 2-9).
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet: duplicate the "2022-Q2" sheet (so it
#    inherits the exact same layout/styles) right before it, rename it,
#    then overwrite its data with the new quarter's figures.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.ActiveSheet
$q3.Name = "2022-Q3"

# Fund-code / figure columns are stored as text in this workbook (leading
# zeros in fund codes, fixed-precision percentage-like strings), so force
# a text number format before assigning them.
foreach ($col in @("B", "D", "E", "F", "G")) {
    $q3.Range($col + "2:" + $col + "5").NumberFormat = "@"
}

# Row 2: 008707 / 建信富时100指数（QDII）美元现汇 A
$q3.Range("B2").Value = "008707"
$q3.Range("C2").Value = "建信富时100指数（QDII）美元现汇 A"
$q3.Range("D2").Value = "0.48"
$q3.Range("E2").Value = "89.38"
$q3.Range("F2").Value = "3.55"
$q3.Range("G2").Value = "0.0170"
$q3.Range("H2").Value = 8

# Row 3: 539003 / 建信富时100指数（QDII）人民币A
$q3.Range("B3").Value = "539003"
$q3.Range("C3").Value = "建信富时100指数（QDII）人民币A"
$q3.Range("D3").Value = "0.48"
$q3.Range("E3").Value = "89.38"
$q3.Range("F3").Value = "3.55"
$q3.Range("G3").Value = "0.0170"
$q3.Range("H3").Value = 8

# Row 4: 008706 / 建信富时100指数（QDII）人民币 C
$q3.Range("B4").Value = "008706"
$q3.Range("C4").Value = "建信富时100指数（QDII）人民币 C"
$q3.Range("D4").Value = "0.19"
$q3.Range("E4").Value = "89.38"
$q3.Range("F4").Value = "3.55"
$q3.Range("G4").Value = "0.0067"
$q3.Range("H4").Value = 8

# Row 5: 008708 / 建信富时100指数（QDII）美元现汇 C
$q3.Range("B5").Value = "008708"
$q3.Range("C5").Value = "建信富时100指数（QDII）美元现汇 C"
$q3.Range("D5").Value = "0.19"
$q3.Range("E5").Value = "89.38"
$q3.Range("F5").Value = "3.55"
$q3.Range("G5").Value = "0.0067"
$q3.Range("H5").Value = 8
